# Fruta / hortaliza, semanal
#
# The published workbook gained one additional weekly price observation.
# A new record is inserted as row 202 (pushing every existing row at/after
# 202 down by one), so the sheet grows from A1:R279 to A1:R280.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 202..279 down to 203..280, opening up a blank row 202
# (formatting of the row above is inherited automatically by Excel).
$ws.Rows.Item(202).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Range("A202").Value = 3
$ws.Range("B202").Value = "Femacal de La Calera"
$ws.Range("C202").Value = "Coquimbo"
$ws.Range("D202").Value = 45134
$ws.Range("E202").Value = 5
$ws.Range("F202").Value = 100112026
$ws.Range("G202").Value = "Haba"
$ws.Range("H202").Value = "Sin especificar"
$ws.Range("I202").Value = "Primera"
$ws.Range("J202").Value = 80
$ws.Range("K202").Value = 15000
$ws.Range("L202").Value = 16000
$ws.Range("M202").Value = 15500
$ws.Range("N202").Value = "$/saco 25 kilos"
$ws.Range("O202").Value = "Provincia de Limarí"
$ws.Range("P202").Value = 620
$ws.Range("Q202").Value = 25
$ws.Range("R202").Value = "Hortaliza"
